$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1286.6666
$ws.Range("J17").Value = 1286.6666
$ws.Range("L17").Value = 3859.9998
$ws.Range("N17").Value = -4195.9998

$ws.Range("H32").Value = 4851.364
$ws.Range("I32").Value = 4594.6
$ws.Range("J32").Value = 5065.3335
$ws.Range("K32").Value = 4594.6
$ws.Range("L32").Value = 5065.3335
$ws.Range("M32").Value = -4268.6
$ws.Range("N32").Value = -5717.3335

$ws.Range("H105").Value = 64000
$ws.Range("J105").Value = 64000
$ws.Range("L105").Value = 64000
$ws.Range("N105").Value = -70988

$ws.Range("H135").Value = 742.61536
$ws.Range("I135").Value = 742.61536
$ws.Range("K135").Value = 6683.53824
$ws.Range("M135").Value = -4148.53824

$ws.Range("H137").Value = 1248.1052
$ws.Range("I137").Value = 1142.6666
$ws.Range("K137").Value = 3427.9998
$ws.Range("M137").Value = -877.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1244.2195
$ws.Range("I32").Value = 1274.6709
$ws.Range("K32").Value = 1274.6709
$ws.Range("M32").Value = -987.6709000000001

$ws.Range("H74").Value = 4950.44
$ws.Range("I74").Value = 5098.9165
$ws.Range("J74").Value = 4813.385
$ws.Range("K74").Value = 5098.9165
$ws.Range("L74").Value = 4813.385
$ws.Range("M74").Value = -4224.9165
$ws.Range("N74").Value = -6561.385

$ws.Range("H77").Value = 4950.44
$ws.Range("I77").Value = 5098.9165
$ws.Range("J77").Value = 4813.385
$ws.Range("K77").Value = 25494.5825
$ws.Range("L77").Value = 24066.925
$ws.Range("M77").Value = -21126.5825
$ws.Range("N77").Value = -32802.925

$ws.Range("H80").Value = 42839.2
$ws.Range("J80").Value = 42839.2
$ws.Range("L80").Value = 42839.2
$ws.Range("N80").Value = -44835.2

$ws.Range("H83").Value = 42839.2
$ws.Range("J83").Value = 42839.2
$ws.Range("L83").Value = 128517.6
$ws.Range("N83").Value = -138501.6

$ws.Range("H110").Value = 6048.7144
$ws.Range("I110").Value = 3962.9285
$ws.Range("K110").Value = 3962.9285
$ws.Range("M110").Value = -1917.9285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 44653.8
$ws.Range("J35").Value = 49367.25
$ws.Range("L35").Value = 49367.25
$ws.Range("N35").Value = -49987.25

$ws.Range("H82").Value = 30800
$ws.Range("I82").Value = 10000
$ws.Range("K82").Value = 10000
$ws.Range("M82").Value = -9617

$ws.Range("H85").Value = 30800
$ws.Range("I85").Value = 10000
$ws.Range("K85").Value = 10000
$ws.Range("M85").Value = -8674

$ws.Range("H99").Value = 4514.5356
$ws.Range("I99").Value = 3205.6316
$ws.Range("K99").Value = 3205.6316
$ws.Range("M99").Value = -1707.6316

$ws.Range("H134").Value = 4027.0732
$ws.Range("I134").Value = 4027.0732
$ws.Range("K134").Value = 12081.2196
$ws.Range("M134").Value = -9546.2196

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 41674.668
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 41674.668
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 41674.668
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -43172.668

$ws.Range("H71").Value = 41674.668
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 41674.668
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 125024.004
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -132512.004

$ws.Range("H74").Value = 41499.332
$ws.Range("J74").Value = 41499.332
$ws.Range("L74").Value = 41499.332
$ws.Range("N74").Value = -43247.332

$ws.Range("H77").Value = 41499.332
$ws.Range("J77").Value = 41499.332
$ws.Range("L77").Value = 124497.996
$ws.Range("N77").Value = -133233.996

$ws.Range("H122").Value = 3103.2693
$ws.Range("I122").Value = 2757.3
$ws.Range("J122").Value = 4256.5
$ws.Range("K122").Value = 8271.900000000001
$ws.Range("L122").Value = 12769.5
$ws.Range("M122").Value = -5821.900000000001
$ws.Range("N122").Value = -17669.5

$ws.Range("H132").Value = 2053.8572
$ws.Range("I132").Value = 1729.5
$ws.Range("K132").Value = 5188.5
$ws.Range("M132").Value = -2658.5

$ws.Range("H135").Value = 81322.555
$ws.Range("J135").Value = 81737.875
$ws.Range("L135").Value = 81737.875
$ws.Range("N135").Value = -91877.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 7214934.5
$ws.Range("I121").Value = 508.5
$ws.Range("J121").Value = 16834170
$ws.Range("K121").Value = 1525.5
$ws.Range("L121").Value = 50502510
$ws.Range("M121").Value = -215.5
$ws.Range("N121").Value = -50505130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 45109.668
$ws.Range("J46").Value = 51552.668
$ws.Range("L46").Value = 51552.668
$ws.Range("N46").Value = -51864.668

$ws.Range("H122").Value = 3767.8333
$ws.Range("I122").Value = 2164.25
$ws.Range("J122").Value = 6975
$ws.Range("K122").Value = 6492.75
$ws.Range("L122").Value = 20925
$ws.Range("M122").Value = -4042.75
$ws.Range("N122").Value = -25825

$ws.Range("H132").Value = 5414.875
$ws.Range("J132").Value = 9833
$ws.Range("L132").Value = 29499
$ws.Range("N132").Value = -34559

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2950.25
$ws.Range("J16").Value = 2900.5
$ws.Range("L16").Value = 2900.5
$ws.Range("N16").Value = -3240.5

$ws.Range("H22").Value = 2249.4285
$ws.Range("I22").Value = 2354.889
$ws.Range("J22").Value = 2059.6
$ws.Range("K22").Value = 2354.889
$ws.Range("L22").Value = 2059.6
$ws.Range("M22").Value = -2059.889
$ws.Range("N22").Value = -2649.6

$ws.Range("H27").Value = 2249.4285
$ws.Range("I27").Value = 2354.889
$ws.Range("J27").Value = 2059.6
$ws.Range("K27").Value = 2354.889
$ws.Range("L27").Value = 2059.6
$ws.Range("M27").Value = -2247.889
$ws.Range("N27").Value = -2273.6

$ws.Range("H46").Value = 10804.6875
$ws.Range("I46").Value = 2527.923
$ws.Range("K46").Value = 2527.923
$ws.Range("M46").Value = -2339.923

$ws.Range("H132").Value = 10010.433
$ws.Range("I132").Value = 11343.556
$ws.Range("K132").Value = 34030.66800000001
$ws.Range("M132").Value = -31500.66800000001

$ws.Range("H136").Value = 8566
$ws.Range("I136").Value = 8099.875
$ws.Range("J136").Value = 9498.25
$ws.Range("K136").Value = 24299.625
$ws.Range("L136").Value = 28494.75
$ws.Range("M136").Value = -21749.625
$ws.Range("N136").Value = -33594.75

$ws.Range("H140").Value = 95850.836
$ws.Range("J140").Value = 95850.836
$ws.Range("L140").Value = 95850.836
$ws.Range("N140").Value = -106210.836
